# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.888.22"
$ws.Range("E2").Value = "  +2.55%  "
$ws.Range("D3").Value = "'1.899.85"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("D5").Value = "'245.86"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").Value = "'1.0000"
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("D7").Value = "'0.5014"
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("D8").Value = "'0.2977"
$ws.Range("E8").Value = "  +1.89%  "
$ws.Range("D9").Value = "'0.06829"
$ws.Range("E9").Value = "  +2.99%  "
$ws.Range("D10").Value = "'1.904.24"
$ws.Range("E10").Value = "  +1.08%  "
$ws.Range("D11").Value = "'17.11"
$ws.Range("E11").Value = "  +1.34%  "
$ws.Range("D12").Value = "'0.07331"
$ws.Range("E12").Value = "  +1.75%  "
$ws.Range("D13").Value = "'91.73"
$ws.Range("E13").Value = "  +7.01%  "
$ws.Range("D14").Value = "'5.096"
$ws.Range("E14").Value = "  +5.36%  "
$ws.Range("D15").Value = "'0.6781"
$ws.Range("E15").Value = "  +1.81%  "
$ws.Range("D16").Value = "'30.856.21"
$ws.Range("E16").Value = "  +2.52%  "
$ws.Range("D17").Value = "'0.000008018"
$ws.Range("E17").Value = "  +2.35%  "
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("E19").Value = "  +3.60%  "
$ws.Range("D20").Value = "'2.150.41"
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").Value = "'4.873"
$ws.Range("E22").Value = "  +2.26%  "
$ws.Range("D23").Value = "'182.21"
$ws.Range("E23").Value = "  +34.66%  "
$ws.Range("D24").Value = "'6.075"
$ws.Range("E24").Value = "  +8.26%  "
$ws.Range("D25").Value = "'9.351"
$ws.Range("E25").Value = "  +1.88%  "
$ws.Range("D26").Value = "'154.85"
$ws.Range("E26").Value = "  +2.46%  "
$ws.Range("D27").Value = "'18.65"
$ws.Range("E27").Value = "  +10.99%  "
$ws.Range("D28").Value = "'1.944"
$ws.Range("E28").Value = "  +1.63%  "
$ws.Range("D29").Value = "'1.394"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("D30").Value = "'4.341"
$ws.Range("E30").Value = "  +4.20%  "
$ws.Range("D31").Value = "'0.08990"
$ws.Range("E31").Value = "  +3.60%  "
$ws.Range("D32").Value = "'4.059"
$ws.Range("D33").Value = "'0.05268"
$ws.Range("E33").Value = "  +5.24%  "
$ws.Range("D34").Value = "'0.7464"
$ws.Range("E34").Value = "  +5.70%  "
$ws.Range("D35").Value = "'1.136"
$ws.Range("E35").Value = "  +2.63%  "
$ws.Range("D36").Value = "'2.669"
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("D37").Value = "'0.01933"
$ws.Range("E37").Value = "  +17.11%  "
$ws.Range("D38").Value = "'2.725"
$ws.Range("E38").Value = "  +0.86%  "
$ws.Range("D39").Value = "'2.184"
$ws.Range("E39").Value = "  -0.73%  "
$ws.Range("D40").Value = "'0.9373"
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("D41").Value = "'0.4392"
$ws.Range("E41").Value = "  +4.43%  "
$ws.Range("D42").Value = "'105.89"
$ws.Range("E42").Value = "  +3.97%  "
$ws.Range("D43").Value = "'5.833"
$ws.Range("E43").Value = "  -2.18%  "
$ws.Range("D44").Value = "'1.001"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D46").Value = "'0.1345"
$ws.Range("D48").Value = "'0.3919"
$ws.Range("E48").Value = "  +5.41%  "
$ws.Range("D49").Value = "'8.578"
$ws.Range("E49").Value = "  +3.66%  "
$ws.Range("D50").Value = "'33.32"
$ws.Range("E50").Value = "  +2.70%  "
$ws.Range("D51").Value = "'1.384"
$ws.Range("E51").Value = "  +2.49%  "
